# Estadisticos Segundo Parcial Sin Ameca
# Adds a new subject/group "Ingles II" / "2ALCV" to the statistics sheets
# and adds the corresponding "Rescatables" (students still needing to pass)
# rows for that new group, while keeping the existing data intact.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheets 1-3: "Estadisticos 1P", "Estadisticos 2P", "Estadisticos Final"
# Each sheet currently has rows 2..5 (row 2 = Formación socioemocional II,
# row 3 = Pensamiento matemático II / 2APV, row 4 = Pensamiento matemático II
# / 2ASV, row 5 = MANTIENE... / 4BEM).
# A brand new row must be inserted right before the current row 3 holding the
# new "Ingles II" / "2ALCV" totals.
# ---------------------------------------------------------------------------

$statSheetNames = @("Estadisticos 1P", "Estadisticos 2P", "Estadisticos Final")

# Values for the new "Ingles II" / "2ALCV" row, per sheet (C,D,E,F,G,H)
$newRowValues = @{
    "Estadisticos 1P"     = @(28, 0, 16, 12, 42.86, 5.7)
    "Estadisticos 2P"     = @(28, 0, 6, 22, 78.57, 5.7)
    "Estadisticos Final"  = @(28, 0, 6, 22, 78.57, 6.8)
}

# Updated H (Promedio) values for the existing "MANTIENE..." / 4BEM row
# (previously row 5, now row 6) per sheet.
$lastRowH = @{
    "Estadisticos 1P"    = 7.3
    "Estadisticos 2P"    = 7.3
    "Estadisticos Final" = 8
}

# Updated H (Promedio) values for the "Pensamiento matemático II" / 2ASV row
# (previously row 4, now row 5) per sheet.
$row5H = @{
    "Estadisticos 1P"    = 6.8
    "Estadisticos 2P"    = 6.8
    "Estadisticos Final" = 7.1
}

foreach ($name in $statSheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Insert a new row above row 3 - everything from the old row 3 downward
    # shifts down by one.
    $ws.Rows.Item(3).Insert()

    $vals = $newRowValues[$name]

    $ws.Range("A3").Value = "Ingles II"
    $ws.Range("B3").Value = "2ALCV"
    $ws.Range("C3").Value = $vals[0]
    $ws.Range("D3").Value = $vals[1]
    $ws.Range("E3").Value = $vals[2]
    $ws.Range("F3").Value = $vals[3]
    $ws.Range("G3").Value = $vals[4]
    $ws.Range("H3").Value = $vals[5]

    # Row 4 (old row 3, "Pensamiento matemático II" / 2APV) keeps its values.
    # Row 5 (old row 4, "Pensamiento matemático II" / 2ASV) - update H.
    $ws.Range("H5").Value = $row5H[$name]

    # Row 6 (old row 5, "MANTIENE..." / 4BEM) - update H.
    $ws.Range("H6").Value = $lastRowH[$name]
}

# ---------------------------------------------------------------------------
# Sheet 4: "Rescatables"
# New rows need to be inserted among the existing ones to keep them grouped
# by subject/group, and a few "Promedio"/G values changed because totals are
# recomputed per-group after adding the "Ingles II" group.
# Final layout (rows 2-12):
#  2: 24330051920330 VASQUEZ   PEREZ     DANIELA LILI      Ingles II                 2ALCV 4
#  3: 24330051920246 ZUNO      FLORES    ALIN MARIEL       Ingles II                 2ALCV 4
#  4: 24330051920393 MUÑOZ     REYES     ERWIN ISRAEL      Pensamiento matemático II 2APV  4
#  5: 24330051920274 CLEMENTE  JUAREZ    BRYAN             Pensamiento matemático II 2APV  3
#  6: 24330051920255 LOPEZ     ROSAS     ERNESTO           Pensamiento matemático II 2APV  3
#  7: 23330051920312 VERA      PEREZ     ALEYDA MONSERRAT  Pensamiento matemático II 2ASV  3
#  8: 24330051920220 GARCIA    CHAPARRO  MAYKA XIMENA      Ingles II                 2ALCV 2
#  9: 24330051920226 LEYVA     HERNANDEZ EUNICE GUADALUPE  Ingles II                 2ALCV 2
# 10: 22330051920007 CARRERA   GARCIA    ANA KAREN         MANTIENE...               4BEM  2
# 11: 24330051920396 MARTINEZ  GONZALEZ  SANTIAGO          Ingles II                 2ALCV 1
# 12: 24330051920369 TORRES    GUTIERREZ JESUS ENRIQUE     Pensamiento matemático II 2APV  1
# ---------------------------------------------------------------------------

$rescatables = $wb.Worksheets.Item("Rescatables")

# Insert two new blank rows after the current row 5 (old data) - these become
# rows 8 and 9 in the final sheet (for GARCIA/CHAPARRO and LEYVA/HERNANDEZ).
$rescatables.Rows.Item(6).Insert()
$rescatables.Rows.Item(7).Insert()

# Insert one new blank row that becomes row 11 (for MARTINEZ/GONZALEZ),
# pushing the old row 7 (TORRES) down to row 12.
$rescatables.Rows.Item(9).Insert()

# Now write the full final data set for rows 2-12 explicitly so every value
# (including the ones that only shifted position) ends up correct.
$data = @(
    @(2,  24330051920330, "VASQUEZ",  "PEREZ",     "DANIELA LILI",     "Ingles II",                 "2ALCV", 4),
    @(3,  24330051920246, "ZUNO",     "FLORES",    "ALIN MARIEL",      "Ingles II",                 "2ALCV", 4),
    @(4,  24330051920393, "MUÑOZ",    "REYES",     "ERWIN ISRAEL",     "Pensamiento matemático II", "2APV",  4),
    @(5,  24330051920274, "CLEMENTE", "JUAREZ",    "BRYAN",            "Pensamiento matemático II", "2APV",  3),
    @(6,  24330051920255, "LOPEZ",    "ROSAS",     "ERNESTO",          "Pensamiento matemático II", "2APV",  3),
    @(7,  23330051920312, "VERA",     "PEREZ",     "ALEYDA MONSERRAT", "Pensamiento matemático II", "2ASV",  3),
    @(8,  24330051920220, "GARCIA",   "CHAPARRO",  "MAYKA XIMENA",     "Ingles II",                 "2ALCV", 2),
    @(9,  24330051920226, "LEYVA",    "HERNANDEZ", "EUNICE GUADALUPE", "Ingles II",                 "2ALCV", 2),
    @(10, 22330051920007, "CARRERA",  "GARCIA",    "ANA KAREN",        "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO", "4BEM", 2),
    @(11, 24330051920396, "MARTINEZ", "GONZALEZ",  "SANTIAGO",         "Ingles II",                 "2ALCV", 1),
    @(12, 24330051920369, "TORRES",   "GUTIERREZ", "JESUS ENRIQUE",    "Pensamiento matemático II", "2APV",  1)
)

foreach ($row in $data) {
    $r = $row[0]
    $rescatables.Range("A$r").Value = $row[1]
    $rescatables.Range("B$r").Value = $row[2]
    $rescatables.Range("C$r").Value = $row[3]
    $rescatables.Range("D$r").Value = $row[4]
    $rescatables.Range("E$r").Value = $row[5]
    $rescatables.Range("F$r").Value = $row[6]
    $rescatables.Range("G$r").Value = $row[7]
}
